$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4713212688202191
$ws.Range("D2").Value = 0.001862681134888078
$ws.Range("E2").Value = 0.09643067662128146
$ws.Range("F2").Value = 0.5085473922686248
$ws.Range("G2").Value = 0.002345227075332708
$ws.Range("I2").Value = 1.067482939672075
$ws.Range("M2").Value = 1.547875453301771
$ws.Range("O2").Value = 1.58785648042624

$ws.Range("B3").Value = 0.4140981522994878
$ws.Range("D3").Value = 0.001637195644702771
$ws.Range("E3").Value = 0.0951427238208673
$ws.Range("F3").Value = 0.4822099989406752
$ws.Range("G3").Value = 0.002349194837688415
$ws.Range("I3").Value = 0.9705780176217047
$ws.Range("M3").Value = 1.355775510223609
$ws.Range("O3").Value = 1.511424509157536

$ws.Range("B4").Value = 0.3788171801890883
$ws.Range("D4").Value = 0.001498812952636541
$ws.Range("E4").Value = 0.09456315371560464
$ws.Range("F4").Value = 0.4665311414210436
$ws.Range("G4").Value = 0.002351759318313542
$ws.Range("I4").Value = 0.9114280091377651
$ws.Range("M4").Value = 1.237352577027451
$ws.Range("O4").Value = 1.466148974409663

$ws.Range("B5").Value = 0.3644039094488107
$ws.Range("D5").Value = 0.001442441067782241
$ws.Range("E5").Value = 0.09437993358347896
$ws.Range("F5").Value = 0.4602643968378999
$ws.Range("G5").Value = 0.00235283672369739
$ws.Range("I5").Value = 0.8874136321980188
$ws.Range("M5").Value = 1.188978449266784
$ws.Range("O5").Value = 1.448110087501249

$ws.Range("B6").Value = 0.3620084432925239
$ws.Range("D6").Value = 0.001433081887332932
$ws.Range("E6").Value = 0.09435270417393937
$ws.Range("F6").Value = 0.459231170390197
$ws.Range("G6").Value = 0.002353017583713723
$ws.Range("I6").Value = 0.8834315351695636
$ws.Range("M6").Value = 1.180939058820002
$ws.Range("O6").Value = 1.44513944745762

$ws.Range("B7").Value = 0.3786229425950296
$ws.Range("D7").Value = 0.00149805261485092
$ws.Range("E7").Value = 0.09456046852540467
$ws.Range("F7").Value = 0.4664461315701232
$ws.Range("G7").Value = 0.002351773717301109
$ws.Range("I7").Value = 0.9111037771852892
$ws.Range("M7").Value = 1.236700651442405
$ws.Range("O7").Value = 1.465904036652972

$ws.Range("B8").Value = 0.4516213585723676
$ws.Range("D8").Value = 0.001784921927969307
$ws.Range("E8").Value = 0.0959426488858881
$ws.Range("F8").Value = 0.4993632655622946
$ws.Range("G8").Value = 0.002346568604467449
$ws.Range("I8").Value = 1.033998450626143
$ws.Range("M8").Value = 1.481739276743127
$ws.Range("O8").Value = 1.561156618889157

$ws.Range("B9").Value = 0.5935937729756802
$ws.Range("D9").Value = 0.002347885896298862
$ws.Range("E9").Value = 0.1003373583811253
$ws.Range("F9").Value = 0.5678832833093566
$ws.Range("G9").Value = 0.002337374159377399
$ws.Range("I9").Value = 1.277713511473053
$ws.Range("M9").Value = 1.958399799918084
$ws.Range("O9").Value = 1.761291367323622

$ws.Range("B10").Value = 0.6971673519410047
$ws.Range("D10").Value = 0.002761660278299161
$ws.Range("E10").Value = 0.1046057654617343
$ws.Range("F10").Value = 0.620737043840478
$ws.Range("G10").Value = 0.002331229521099414
$ws.Range("I10").Value = 1.45837596672385
$ws.Range("M10").Value = 2.306143383950541
$ws.Range("O10").Value = 1.916787550815457

$ws.Range("B11").Value = 0.7441242489681485
$ws.Range("D11").Value = 0.002949923056540626
$ws.Range("E11").Value = 0.1067763425217478
$ws.Range("F11").Value = 0.6453473143921968
$ws.Range("G11").Value = 0.002328565262924213
$ws.Range("I11").Value = 1.540904482121022
$ws.Range("M11").Value = 2.463788356980302
$ws.Range("O11").Value = 1.989433752520029

$ws.Range("B12").Value = 0.7618823371220174
$ws.Range("D12").Value = 0.003021216932971527
$ws.Range("E12").Value = 0.1076314372214568
$ws.Range("F12").Value = 0.6547496273801556
$ws.Range("G12").Value = 0.002327575098844085
$ws.Range("I12").Value = 1.572204416953383
$ws.Range("M12").Value = 2.523403785311615
$ws.Range("O12").Value = 2.017223082984287

$ws.Range("B13").Value = 0.7580588676837579
$ws.Range("D13").Value = 0.003005862425261796
$ws.Range("E13").Value = 0.1074457993199331
$ws.Range("F13").Value = 0.6527209595395078
$ws.Range("G13").Value = 0.002327787516698064
$ws.Range("I13").Value = 1.565461293395742
$ws.Range("M13").Value = 2.510568197024668
$ws.Range("O13").Value = 2.011225628075351

$ws.Range("B14").Value = 0.7455856910940497
$ws.Range("D14").Value = 0.002955788401301618
$ws.Range("E14").Value = 0.1068460260331037
$ws.Range("F14").Value = 0.646119177935006
$ws.Range("G14").Value = 0.002328483426732241
$ws.Range("I14").Value = 1.543478588321847
$ws.Range("M14").Value = 2.468694599746811
$ws.Range("O14").Value = 1.991714360693607

$ws.Range("B15").Value = 0.7379424374244081
$ws.Range("D15").Value = 0.002925116927421811
$ws.Range("E15").Value = 0.1064829710947706
$ws.Range("F15").Value = 0.6420862395440707
$ws.Range("G15").Value = 0.00232891212822095
$ws.Range("I15").Value = 1.530019783037915
$ws.Range("M15").Value = 2.443035140440145
$ws.Range("O15").Value = 1.97979973506051

$ws.Range("B16").Value = 0.6940953544213357
$ws.Range("D16").Value = 0.002749357371222061
$ws.Range("E16").Value = 0.1044685381705079
$ws.Range("F16").Value = 0.6191402327384452
$ws.Range("G16").Value = 0.002331406263524889
$ws.Range("I16").Value = 1.452989345951465
$ws.Range("M16").Value = 2.295829692187453
$ws.Range("O16").Value = 1.912078824644368

$ws.Range("B17").Value = 0.6671553873597418
$ws.Range("D17").Value = 0.002641541843235728
$ws.Range("E17").Value = 0.103291528352834
$ws.Range("F17").Value = 0.6052097929921985
$ws.Range("G17").Value = 0.002332969807099393
$ws.Range("I17").Value = 1.405820866799274
$ws.Range("M17").Value = 2.20538216493091
$ws.Range("O17").Value = 1.871027112345701

$ws.Range("B18").Value = 0.6516452742277181
$ws.Range("D18").Value = 0.002579532580831767
$ws.Range("E18").Value = 0.1026360698584554
$ws.Range("F18").Value = 0.5972506393104737
$ws.Range("G18").Value = 0.002333881449130755
$ws.Range("I18").Value = 1.378723281700957
$ws.Range("M18").Value = 2.153307977228252
$ws.Range("O18").Value = 1.847594744034325

$ws.Range("B19").Value = 0.6463912655031265
$ws.Range("D19").Value = 0.002558537960737794
$ws.Range("E19").Value = 0.1024178330734316
$ws.Range("F19").Value = 0.5945649118753238
$ws.Range("G19").Value = 0.002334192236870523
$ws.Range("I19").Value = 1.369554124832575
$ws.Range("M19").Value = 2.135667854263971
$ws.Range("O19").Value = 1.83969161940206

$ws.Range("B20").Value = 0.6700247450202426
$ws.Range("D20").Value = 0.002653018656982198
$ws.Range("E20").Value = 0.1034145935009221
$ws.Range("F20").Value = 0.6066871888612866
$ws.Range("G20").Value = 0.002332802089239892
$ws.Range("I20").Value = 1.410838680635095
$ws.Range("M20").Value = 2.215015774858841
$ws.Range("O20").Value = 1.875378523730603

$ws.Range("B21").Value = 0.7492500067679089
$ws.Range("D21").Value = 0.002970496293382752
$ws.Range("E21").Value = 0.1070212925339931
$ws.Range("F21").Value = 0.6480560188925608
$ws.Range("G21").Value = 0.002328278513700656
$ws.Range("I21").Value = 1.54993414426562
$ws.Range("M21").Value = 2.480996119331792
$ws.Range("O21").Value = 1.997437658612398

$ws.Range("B22").Value = 0.8008912528991345
$ws.Range("D22").Value = 0.003178002116694501
$ws.Range("E22").Value = 0.1095717857874163
$ws.Range("F22").Value = 0.6755770338646556
$ws.Range("G22").Value = 0.002325431235317889
$ws.Range("I22").Value = 1.641121334320758
$ws.Range("M22").Value = 2.65435448866171
$ws.Range("O22").Value = 2.078843583671414

$ws.Range("B23").Value = 0.7733420679550136
$ws.Range("D23").Value = 0.003067251530165294
$ws.Range("E23").Value = 0.1081927715060402
$ws.Range("F23").Value = 0.6608437815068413
$ws.Range("G23").Value = 0.002326940928146354
$ws.Range("I23").Value = 1.592427784857506
$ws.Range("M23").Value = 2.561874270899523
$ws.Range("O23").Value = 2.035244574471449

$ws.Range("B24").Value = 0.6687275768046561
$ws.Range("D24").Value = 0.002647830065903634
$ws.Range("E24").Value = 0.1033588896620614
$ws.Range("F24").Value = 0.6060191036701212
$ws.Range("G24").Value = 0.002332877874887201
$ws.Range("I24").Value = 1.408570063974196
$ws.Range("M24").Value = 2.210660652281007
$ws.Range("O24").Value = 1.873410725215365

$ws.Range("B25").Value = 0.555314266893447
$ws.Range("D25").Value = 0.002195555114017367
$ws.Range("E25").Value = 0.09896687188386721
$ws.Range("F25").Value = 0.5489123345535347
$ws.Range("G25").Value = 0.002339753790606802
$ws.Range("I25").Value = 1.211496843730842
$ws.Range("M25").Value = 1.829874998367501
$ws.Range("O25").Value = 1.705687799287432
